$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-26 21:18:39"
$ws.Range("E3").Value = "2026-02-26 21:18:41"
$ws.Range("L3").Value = "19.4 km/h - 120º 20:42 TU"
$ws.Range("E4").Value = "2026-02-26 21:18:44"
$ws.Range("O4").Value = "10.7 °C"
$ws.Range("E5").Value = "2026-02-26 21:18:46"
$ws.Range("H5").Value = "45%"
$ws.Range("E6").Value = "2026-02-26 21:18:49"
$ws.Range("H6").Value = "84%"
$ws.Range("O6").Value = "11.8 °C"
$ws.Range("E7").Value = "2026-02-26 21:18:52"
$ws.Range("H7").Value = "82%"
$ws.Range("E8").Value = "2026-02-26 21:18:55"
$ws.Range("H8").Value = "89%"
$ws.Range("E9").Value = "2026-02-26 21:18:57"
$ws.Range("N9").Value = "7.2 °C 20:30 TU"
$ws.Range("O9").Value = "12.1 °C"
$ws.Range("E10").Value = "2026-02-26 21:19:00"
$ws.Range("O10").Value = "9.5 °C"
$ws.Range("E11").Value = "2026-02-26 21:19:03"
$ws.Range("E12").Value = "2026-02-26 21:19:05"
$ws.Range("H12").Value = "92%"
$ws.Range("O12").Value = "11.3 °C"
$ws.Range("E13").Value = "2026-02-26 21:19:08"
$ws.Range("J13").Value = "1028.2 hPa"
$ws.Range("E14").Value = "2026-02-26 21:19:11"
$ws.Range("N14").Value = "5.7 °C 20:47 TU"
$ws.Range("O14").Value = "11.4 °C"
$ws.Range("E15").Value = "2026-02-26 21:19:13"
$ws.Range("E16").Value = "2026-02-26 21:19:15"
$ws.Range("E17").Value = "2026-02-26 21:19:18"
$ws.Range("E18").Value = "2026-02-26 21:19:21"
$ws.Range("J18").Value = "1027.3 hPa"
$ws.Range("O18").Value = "11.9 °C"
$ws.Range("E19").Value = "2026-02-26 21:19:23"
$ws.Range("O19").Value = "11.4 °C"
$ws.Range("E20").Value = "2026-02-26 21:19:26"
$ws.Range("H20").Value = "48%"
$ws.Range("O20").Value = "2.6 °C"
$ws.Range("E21").Value = "2026-02-26 21:19:29"
$ws.Range("J21").Value = "1027.0 hPa"
$ws.Range("E22").Value = "2026-02-26 21:19:31"
$ws.Range("E23").Value = "2026-02-26 21:19:34"
$ws.Range("O23").Value = "3.4 °C"
$ws.Range("E24").Value = "2026-02-26 21:19:37"
$ws.Range("O24").Value = "10.4 °C"
$ws.Range("E25").Value = "2026-02-26 21:19:39"
$ws.Range("O25").Value = "5.2 °C"
$ws.Range("E26").Value = "2026-02-26 21:19:42"
$ws.Range("O26").Value = "10.8 °C"
$ws.Range("E27").Value = "2026-02-26 21:19:44"
$ws.Range("E28").Value = "2026-02-26 21:19:47"
$ws.Range("J28").Value = "1026.8 hPa"
$ws.Range("O28").Value = "10.8 °C"
$ws.Range("E29").Value = "2026-02-26 21:19:49"
$ws.Range("N29").Value = "7.2 °C 20:58 TU"
$ws.Range("O29").Value = "11.6 °C"
$ws.Range("E30").Value = "2026-02-26 21:19:52"
$ws.Range("E31").Value = "2026-02-26 21:19:55"
$ws.Range("E32").Value = "2026-02-26 21:19:57"
$ws.Range("H32").Value = "66%"
$ws.Range("O32").Value = "8.0 °C"
$ws.Range("E33").Value = "2026-02-26 21:20:00"
$ws.Range("E34").Value = "2026-02-26 21:20:03"
$ws.Range("E35").Value = "2026-02-26 21:20:05"
$ws.Range("H35").Value = "43%"
$ws.Range("J35").Value = "1025.5 hPa"
$ws.Range("O35").Value = "12.1 °C"
$ws.Range("E36").Value = "2026-02-26 21:20:08"
$ws.Range("O36").Value = "12.5 °C"
$ws.Range("E37").Value = "2026-02-26 21:20:11"
$ws.Range("H37").Value = "74%"
$ws.Range("O37").Value = "7.8 °C"
$ws.Range("E38").Value = "2026-02-26 21:20:13"
$ws.Range("O38").Value = "11.2 °C"
$ws.Range("E39").Value = "2026-02-26 21:20:16"
$ws.Range("E40").Value = "2026-02-26 21:20:19"
$ws.Range("H40").Value = "67%"
$ws.Range("J40").Value = "1027.4 hPa"
$ws.Range("O40").Value = "9.6 °C"
$ws.Range("E41").Value = "2026-02-26 21:20:21"
$ws.Range("O41").Value = "11.1 °C"
$ws.Range("E42").Value = "2026-02-26 21:20:23"
$ws.Range("H42").Value = "88%"
$ws.Range("N42").Value = "7.6 °C 20:59 TU"
$ws.Range("O42").Value = "11.3 °C"
$ws.Range("E43").Value = "2026-02-26 21:20:26"
$ws.Range("H43").Value = "73%"
$ws.Range("K43").Value = "15.7 MJ/m2"
$ws.Range("E44").Value = "2026-02-26 21:20:28"
$ws.Range("H44").Value = "54%"
$ws.Range("E45").Value = "2026-02-26 21:20:31"
$ws.Range("E46").Value = "2026-02-26 21:20:34"
